$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows into the data table (27 -> 30 rows) ---
# This pushes old row 42 (and the two footer rows 47/48) down by 3 rows.
$ws.Rows("42:44").Insert()

# Copy formatting from row 41 (a standard data row) into the 3 new rows
$ws.Range("B41:J41").Copy($ws.Range("B42:J42"))
$ws.Range("B41:J41").Copy($ws.Range("B43:J43"))
$ws.Range("B41:J41").Copy($ws.Range("B44:J44"))

# --- Update summary header figures ---
$ws.Range("E11").Value2 = 1802565    # VALOR MORA total
$ws.Range("C13").Value2 = 5           # Cant. Trabajadores
$ws.Range("F13").Value2 = 26          # Cant. Periodos

# --- Rewrite the full worker/period data table (rows 16-45) ---
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "73186208"
$ws.Range("D16").Value2 = "SEGUNDO JAVIER SALGADO REYES"
$ws.Range("E16").Value2 = "2507"
$ws.Range("F16").Value2 = 2467
$ws.Range("G16").Value2 = 828116
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "9299898"
$ws.Range("D17").Value2 = "FRESMAN MARTINEZ PANZA"
$ws.Range("E17").Value2 = "2502"
$ws.Range("F17").Value2 = 1898
$ws.Range("G17").Value2 = 1423500
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1002388835"
$ws.Range("D18").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E18").Value2 = "2507"
$ws.Range("F18").Value2 = 66000
$ws.Range("G18").Value2 = 1650000
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1002388835"
$ws.Range("D19").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E19").Value2 = "2506"
$ws.Range("F19").Value2 = 66000
$ws.Range("G19").Value2 = 1650000
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1002388835"
$ws.Range("D20").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E20").Value2 = "2505"
$ws.Range("F20").Value2 = 66000
$ws.Range("G20").Value2 = 1650000
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1002388835"
$ws.Range("D21").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E21").Value2 = "2504"
$ws.Range("F21").Value2 = 66000
$ws.Range("G21").Value2 = 1650000
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1002388835"
$ws.Range("D22").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E22").Value2 = "2503"
$ws.Range("F22").Value2 = 66000
$ws.Range("G22").Value2 = 1650000
$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "1002388835"
$ws.Range("D23").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E23").Value2 = "2502"
$ws.Range("F23").Value2 = 66000
$ws.Range("G23").Value2 = 1650000
$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "1002388835"
$ws.Range("D24").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E24").Value2 = "2501"
$ws.Range("F24").Value2 = 66000
$ws.Range("G24").Value2 = 1650000
$ws.Range("B25").Value2 = "CC"
$ws.Range("C25").Value2 = "1002388835"
$ws.Range("D25").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E25").Value2 = "2412"
$ws.Range("F25").Value2 = 66000
$ws.Range("G25").Value2 = 1650000
$ws.Range("B26").Value2 = "CC"
$ws.Range("C26").Value2 = "1002388835"
$ws.Range("D26").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E26").Value2 = "2411"
$ws.Range("F26").Value2 = 66000
$ws.Range("G26").Value2 = 1650000
$ws.Range("B27").Value2 = "CC"
$ws.Range("C27").Value2 = "1002388835"
$ws.Range("D27").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E27").Value2 = "2410"
$ws.Range("F27").Value2 = 66000
$ws.Range("G27").Value2 = 1650000
$ws.Range("B28").Value2 = "CC"
$ws.Range("C28").Value2 = "1002388835"
$ws.Range("D28").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E28").Value2 = "2409"
$ws.Range("F28").Value2 = 66000
$ws.Range("G28").Value2 = 1650000
$ws.Range("B29").Value2 = "CC"
$ws.Range("C29").Value2 = "1002388835"
$ws.Range("D29").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E29").Value2 = "2408"
$ws.Range("F29").Value2 = 66000
$ws.Range("G29").Value2 = 1650000
$ws.Range("B30").Value2 = "CC"
$ws.Range("C30").Value2 = "1002388835"
$ws.Range("D30").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E30").Value2 = "2407"
$ws.Range("F30").Value2 = 66000
$ws.Range("G30").Value2 = 1650000
$ws.Range("B31").Value2 = "CC"
$ws.Range("C31").Value2 = "1002388835"
$ws.Range("D31").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E31").Value2 = "2406"
$ws.Range("F31").Value2 = 66000
$ws.Range("G31").Value2 = 1650000
$ws.Range("B32").Value2 = "CC"
$ws.Range("C32").Value2 = "1002388835"
$ws.Range("D32").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E32").Value2 = "2405"
$ws.Range("F32").Value2 = 66000
$ws.Range("G32").Value2 = 1650000
$ws.Range("B33").Value2 = "CC"
$ws.Range("C33").Value2 = "1002388835"
$ws.Range("D33").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E33").Value2 = "2404"
$ws.Range("F33").Value2 = 66000
$ws.Range("G33").Value2 = 1650000
$ws.Range("B34").Value2 = "CC"
$ws.Range("C34").Value2 = "1002388835"
$ws.Range("D34").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E34").Value2 = "2403"
$ws.Range("F34").Value2 = 66000
$ws.Range("G34").Value2 = 1650000
$ws.Range("B35").Value2 = "CC"
$ws.Range("C35").Value2 = "1002388835"
$ws.Range("D35").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E35").Value2 = "2402"
$ws.Range("F35").Value2 = 66000
$ws.Range("G35").Value2 = 1650000
$ws.Range("B36").Value2 = "CC"
$ws.Range("C36").Value2 = "1002388835"
$ws.Range("D36").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E36").Value2 = "2401"
$ws.Range("F36").Value2 = 66000
$ws.Range("G36").Value2 = 1650000
$ws.Range("B37").Value2 = "CC"
$ws.Range("C37").Value2 = "1002388835"
$ws.Range("D37").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E37").Value2 = "2312"
$ws.Range("F37").Value2 = 66000
$ws.Range("G37").Value2 = 1650000
$ws.Range("B38").Value2 = "CC"
$ws.Range("C38").Value2 = "1002388835"
$ws.Range("D38").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E38").Value2 = "2311"
$ws.Range("F38").Value2 = 66000
$ws.Range("G38").Value2 = 1650000
$ws.Range("B39").Value2 = "CC"
$ws.Range("C39").Value2 = "1002388835"
$ws.Range("D39").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E39").Value2 = "2310"
$ws.Range("F39").Value2 = 66000
$ws.Range("G39").Value2 = 1650000
$ws.Range("B40").Value2 = "CC"
$ws.Range("C40").Value2 = "1002388835"
$ws.Range("D40").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E40").Value2 = "2309"
$ws.Range("F40").Value2 = 66000
$ws.Range("G40").Value2 = 1650000
$ws.Range("B41").Value2 = "CC"
$ws.Range("C41").Value2 = "1002388835"
$ws.Range("D41").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E41").Value2 = "2308"
$ws.Range("F41").Value2 = 66000
$ws.Range("G41").Value2 = 1650000
$ws.Range("B42").Value2 = "CC"
$ws.Range("C42").Value2 = "1002388835"
$ws.Range("D42").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E42").Value2 = "2306"
$ws.Range("F42").Value2 = 66000
$ws.Range("G42").Value2 = 1650000
$ws.Range("B43").Value2 = "CC"
$ws.Range("C43").Value2 = "1002388835"
$ws.Range("D43").Value2 = "JERSON AMAURY PIEDRAHITA QUINTANA"
$ws.Range("E43").Value2 = "2305"
$ws.Range("F43").Value2 = 66000
$ws.Range("G43").Value2 = 1650000
$ws.Range("B44").Value2 = "CC"
$ws.Range("C44").Value2 = "20506516"
$ws.Range("D44").Value2 = "MERLY CHAVEZ NIGRINIS"
$ws.Range("E44").Value2 = "2407"
$ws.Range("F44").Value2 = 60000
$ws.Range("G44").Value2 = 1500000
$ws.Range("B45").Value2 = "CC"
$ws.Range("C45").Value2 = "1143339954"
$ws.Range("D45").Value2 = "CHRISTOPHER RAMIREZ YEPES"
$ws.Range("E45").Value2 = "2507"
$ws.Range("F45").Value2 = 22200
$ws.Range("G45").Value2 = 1850000
